# "Generate Report for Handback"
#
# This fills in the "Latest Target File" / "Latest Handback File" columns
# (F/G) on the per-language sheets, updates the handback status + handback
# datetime, and marks the overall status as handed back.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# Hyperlink font, matching the workbook's existing custom "HyperLink" style
# (underline + RGB 6495ED, stored as a BGR-ordered OLE color value).
$linkUnderline = 2
$linkColor = 15570276

function Set-LinkCell {
    param(
        $ws,
        [string]$cellRef,
        [string]$text,
        [string]$url
    )

    $rng = $ws.Range($cellRef)
    $rng.Value = $text
    $ws.Hyperlinks.Add($rng, $url, "", "", $text) | Out-Null
    $rng.Font.Underline = $linkUnderline
    $rng.Font.Color = $linkColor
}

# ---------------------------------------------------------------------
# Overview sheet: Status column reflects the new handback state as well.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

Set-LinkCell $wsZh "F2" "a136be05-d50e-46d1-9f44-17a045db5ffb.md" "https://github.com/OpenLocalizationTest/oltest/blob/5ff0fcd8a38616703c67fbeb671759135d0859b7/e2e/a136be05-d50e-46d1-9f44-17a045db5ffb.md"
Set-LinkCell $wsZh "G2" "a136be05-d50e-46d1-9f44-17a045db5ffb.638e5344383a04a5ca06e7e6665c050d0eccfaa7.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b623e96f1cdb12910346c94824398fe59962f133/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a136be05-d50e-46d1-9f44-17a045db5ffb.638e5344383a04a5ca06e7e6665c050d0eccfaa7.zh-cn.xlf"

Set-LinkCell $wsZh "F3" "a9227bf9-d672-45db-85a6-1fe97592d078.md" "https://github.com/OpenLocalizationTest/oltest/blob/5ff0fcd8a38616703c67fbeb671759135d0859b7/e2e/a9227bf9-d672-45db-85a6-1fe97592d078.md"
Set-LinkCell $wsZh "G3" "a9227bf9-d672-45db-85a6-1fe97592d078.ba7033e31937687134f690efee6fd988aab332ff.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b623e96f1cdb12910346c94824398fe59962f133/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a9227bf9-d672-45db-85a6-1fe97592d078.ba7033e31937687134f690efee6fd988aab332ff.zh-cn.xlf"

# zh-cn handback completed at 04:37:40 -- replaces the "0001-01-01" placeholder.
$wsZh.Range("H2").Value = "2016-03-20 04:37:40"
$wsZh.Range("H3").Value = "2016-03-20 04:37:40"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

Set-LinkCell $wsDe "F2" "a136be05-d50e-46d1-9f44-17a045db5ffb.md" "https://github.com/OpenLocalizationTest/oltest/blob/5ff0fcd8a38616703c67fbeb671759135d0859b7/e2e/a136be05-d50e-46d1-9f44-17a045db5ffb.md"
Set-LinkCell $wsDe "G2" "a136be05-d50e-46d1-9f44-17a045db5ffb.638e5344383a04a5ca06e7e6665c050d0eccfaa7.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3c0df69927bed3d188868d85e737bd5e86a695f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a136be05-d50e-46d1-9f44-17a045db5ffb.638e5344383a04a5ca06e7e6665c050d0eccfaa7.de-de.xlf"

Set-LinkCell $wsDe "F3" "a9227bf9-d672-45db-85a6-1fe97592d078.md" "https://github.com/OpenLocalizationTest/oltest/blob/5ff0fcd8a38616703c67fbeb671759135d0859b7/e2e/a9227bf9-d672-45db-85a6-1fe97592d078.md"
Set-LinkCell $wsDe "G3" "a9227bf9-d672-45db-85a6-1fe97592d078.ba7033e31937687134f690efee6fd988aab332ff.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e3c0df69927bed3d188868d85e737bd5e86a695f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a9227bf9-d672-45db-85a6-1fe97592d078.ba7033e31937687134f690efee6fd988aab332ff.de-de.xlf"

# de-de handback completed at 04:37:46 (distinct from zh-cn's timestamp).
$wsDe.Range("H2").Value = "2016-03-20 04:37:46"
$wsDe.Range("H3").Value = "2016-03-20 04:37:46"
